$wb = $excel.ActiveWorkbook
foreach ($ws in $wb.Worksheets) {
    $ws.Name = $ws.Name + "_1"
}
foreach ($ws in $wb.Worksheets) {
    Write-Output $ws.Name
}
